# This script applies a row-data permutation to rows 6,7,8,9,11,12 of the
# active worksheet. The columns A,B,D,E,F,G,H,Q,R are shuffled between these
# rows (two independent cycles: 6->9->7->12->6 and 8->11->8), while all other
# columns (C,I,P,S,T,U,V,W,Y,Z,AA,AB,AD,AE,AG,AT,AW,AX,AY) remain untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for each row after the edit (taken directly from the diff).
$data = @{
    6  = @{ A = 111936776; B = 77515; D = "NT"; E = 6425;  F = "Garnlav";     G = "Alectoria sarmentosa";          H = "(Ach.) Ach.";                          Q = 490397.9705776197; R = 7088444.768114219 }
    7  = @{ A = 111936768; B = 90087; D = "LC"; E = 3298;  F = "Trådticka"; G = "Climacocystis borealis";       H = "(Fr.) Kotl. & Pouzar";                 Q = 490316.6991760527; R = 7088522.021911296 }
    8  = @{ A = 111936767; B = 90087; D = "LC"; E = 3298;  F = "Trådticka"; G = "Climacocystis borealis";       H = "(Fr.) Kotl. & Pouzar";                 Q = 490377.1611957431; R = 7088411.830052498 }
    9  = @{ A = 111936777; B = 77515; D = "NT"; E = 6425;  F = "Garnlav";     G = "Alectoria sarmentosa";          H = "(Ach.) Ach.";                          Q = 490055.5835512968; R = 7088708.521274347 }
    11 = @{ A = 111936775; B = 89419; D = "NT"; E = 1204;  F = "Gränsticka"; G = "Phellopilus nigrolimitatus"; H = "(Romell) Niemelä, T.Wagner & M.Fisch."; Q = 490379.7033068824; R = 7088378.819891299 }
    12 = @{ A = 111936781; B = 89793; D = "LC"; E = 4217;  F = "Blodticka";   G = "Meruliopsis taxicola";          H = "(Pers.:Fr.) Bondartsev";               Q = 490314.5888938977; R = 7088551.949221384 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("A$row").Value = $vals.A
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("F$row").Value = $vals.F
    $ws.Range("G$row").Value = $vals.G
    $ws.Range("H$row").Value = $vals.H
    $ws.Range("Q$row").Value = $vals.Q
    $ws.Range("R$row").Value = $vals.R
}
